$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.244.94'
$ws.Range('E2').Value = '  +6.32%  '
$ws.Range('D3').Value = '3.676.72'
$ws.Range('E3').Value = '  +18.36%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Value = '''603.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.56%  '
$ws.Range('D6').Value = '''182.15'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.37%  '
$ws.Range('D7').Value = '3.674.02'
$ws.Range('E7').Value = '  +18.39%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '''0.537'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.87%  '
$ws.Range('E10').Value = '  +7.66%  '
$ws.Range('D11').Value = '''6.62'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.19%  '
$ws.Range('E12').Value = '  +6.76%  '
$ws.Range('D13').Value = '''40.65'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +12.26%  '
$ws.Range('E14').Value = '  +5.62%  '
$ws.Range('D15').Value = '4.286.15'
$ws.Range('E15').Value = '  +18.27%  '
$ws.Range('D16').Value = '71.204.32'
$ws.Range('E16').Value = '  +6.34%  '
$ws.Range('D17').Value = '3.672.06'
$ws.Range('E17').Value = '  +18.17%  '
$ws.Range('E18').Value = '  +1.16%  '
$ws.Range('D19').Value = '''7.51'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.37%  '
$ws.Range('E20').Value = '  +0.87%  '
$ws.Range('D21').Value = '''520.39'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.26%  '
$ws.Range('D22').Value = '''9.25'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +19.33%  '
$ws.Range('D23').Value = '''0.746'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.05%  '
$ws.Range('D24').Value = '''88.29'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.51%  '
$ws.Range('D25').Value = '''2.48'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.94%  '
$ws.Range('D26').Value = '''13.54'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.20%  '
$ws.Range('D27').Value = '''11.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.30%  '
$ws.Range('D28').Value = '''0.999'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').Value = '''2.55'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +11.76%  '
$ws.Range('D30').Value = '''8.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.03%  '
$ws.Range('D31').Value = '''2.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.69%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').Value = '''0.0000111'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +18.53%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '''31.82'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +13.47%  '
$ws.Range('D34').Value = '''0.116'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.03%  '
$ws.Range('D35').Value = '''1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').Value = '''6.14'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.61%  '
$ws.Range('E37').Value = '  +8.03%  '
$ws.Range('D38').Value = '''0.346'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +12.05%  '
$ws.Range('E39').Value = '  +10.24%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '''51.43'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.60%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.130'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.61%  '
$ws.Range('D42').Value = '''45.28'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.27%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '3.135.38'
$ws.Range('E43').Value = '  +12.25%  '
$ws.Range('B44').Value = 'Cosmos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D44').Value = '''8.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.37%  '
$ws.Range('D45').Value = '''415.51'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +11.63%  '
$ws.Range('D46').Value = '''2.78'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.80%  '
$ws.Range('E47').Value = '  +7.05%  '
$ws.Range('D48').Value = '''28.47'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +14.13%  '
$ws.Range('D49').Value = '''139.81'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.35%  '
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').Value = '''2.47'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.73%  '
